$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 499
$ws.Cells.Item(2, 9).Value = 499
$ws.Cells.Item(2, 11).Value = 499
$ws.Cells.Item(2, 13).Value = -386
$ws.Cells.Item(33, 8).Value = 2939.6924
$ws.Cells.Item(33, 9).Value = 499.33334
$ws.Cells.Item(33, 10).Value = 3671.8
$ws.Cells.Item(33, 11).Value = 499.33334
$ws.Cells.Item(33, 12).Value = 3671.8
$ws.Cells.Item(33, 13).Value = -270.33334
$ws.Cells.Item(33, 14).Value = -4129.8
$ws.Cells.Item(40, 8).Value = 3549.1035
$ws.Cells.Item(40, 9).Value = 1769.75
$ws.Cells.Item(40, 11).Value = 1769.75
$ws.Cells.Item(40, 13).Value = -1594.75
$ws.Cells.Item(62, 8).Value = 26793.03
$ws.Cells.Item(62, 9).Value = 29853.654
$ws.Cells.Item(62, 11).Value = 29853.654
$ws.Cells.Item(62, 13).Value = -29229.654
$ws.Cells.Item(65, 8).Value = 26793.03
$ws.Cells.Item(65, 9).Value = 29853.654
$ws.Cells.Item(65, 11).Value = 149268.27
$ws.Cells.Item(65, 13).Value = -146148.27
$ws.Cells.Item(80, 8).Value = 2323.8518
$ws.Cells.Item(80, 9).Value = 619.63635
$ws.Cells.Item(80, 10).Value = 3495.5
$ws.Cells.Item(80, 11).Value = 1858.90905
$ws.Cells.Item(80, 12).Value = 10486.5
$ws.Cells.Item(80, 13).Value = -860.90905
$ws.Cells.Item(80, 14).Value = -12482.5
$ws.Cells.Item(83, 8).Value = 2323.8518
$ws.Cells.Item(83, 9).Value = 619.63635
$ws.Cells.Item(83, 10).Value = 3495.5
$ws.Cells.Item(83, 11).Value = 5576.72715
$ws.Cells.Item(83, 12).Value = 31459.5
$ws.Cells.Item(83, 13).Value = -584.7271499999997
$ws.Cells.Item(83, 14).Value = -41443.5
$ws.Cells.Item(112, 8).Value = 1400.5625
$ws.Cells.Item(112, 10).Value = 1407.2667
$ws.Cells.Item(112, 12).Value = 4221.800099999999
$ws.Cells.Item(112, 14).Value = -6437.800099999999

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 475
$ws.Cells.Item(4, 9).Value = 450
$ws.Cells.Item(4, 11).Value = 450
$ws.Cells.Item(4, 13).Value = -334
$ws.Cells.Item(5, 8).Value = 234.16667
$ws.Cells.Item(5, 10).Value = 351.25
$ws.Cells.Item(5, 12).Value = 351.25
$ws.Cells.Item(5, 14).Value = -575.25
$ws.Cells.Item(32, 8).Value = 4927.6416
$ws.Cells.Item(32, 9).Value = 3610.7869
$ws.Cells.Item(32, 11).Value = 3610.7869
$ws.Cells.Item(32, 13).Value = -3323.7869
$ws.Cells.Item(45, 8).Value = 17428.572
$ws.Cells.Item(45, 9).Value = 34666.668
$ws.Cells.Item(45, 10).Value = 4500
$ws.Cells.Item(45, 11).Value = 34666.668
$ws.Cells.Item(45, 12).Value = 4500
$ws.Cells.Item(45, 13).Value = -34289.668
$ws.Cells.Item(45, 14).Value = -5254
$ws.Cells.Item(61, 8).Value = 5134.1333
$ws.Cells.Item(61, 9).Value = 4786.643
$ws.Cells.Item(61, 10).Value = 9999
$ws.Cells.Item(61, 11).Value = 4786.643
$ws.Cells.Item(61, 12).Value = 9999
$ws.Cells.Item(61, 13).Value = -4574.643
$ws.Cells.Item(61, 14).Value = -10423
$ws.Cells.Item(80, 8).Value = 500024670
$ws.Cells.Item(80, 10).Value = 600022000
$ws.Cells.Item(80, 12).Value = 600022000
$ws.Cells.Item(80, 14).Value = -600023996
$ws.Cells.Item(83, 8).Value = 500024670
$ws.Cells.Item(83, 10).Value = 600022000
$ws.Cells.Item(83, 12).Value = 1800066000
$ws.Cells.Item(83, 14).Value = -1800075984
$ws.Cells.Item(88, 8).Value = 2501.5625
$ws.Cells.Item(88, 9).Value = 2503.875
$ws.Cells.Item(88, 10).Value = 2499.25
$ws.Cells.Item(88, 11).Value = 2503.875
$ws.Cells.Item(88, 12).Value = 2499.25
$ws.Cells.Item(88, 13).Value = -2097.875
$ws.Cells.Item(88, 14).Value = -3311.25
$ws.Cells.Item(91, 8).Value = 2501.5625
$ws.Cells.Item(91, 9).Value = 2503.875
$ws.Cells.Item(91, 10).Value = 2499.25
$ws.Cells.Item(91, 11).Value = 2503.875
$ws.Cells.Item(91, 12).Value = 2499.25
$ws.Cells.Item(91, 13).Value = -1099.875
$ws.Cells.Item(91, 14).Value = -5307.25
$ws.Cells.Item(132, 8).Value = 11353.223
$ws.Cells.Item(132, 9).Value = 10311.357
$ws.Cells.Item(132, 11).Value = 30934.071
$ws.Cells.Item(132, 13).Value = -28404.071
$ws.Cells.Item(136, 8).Value = 5134.1333
$ws.Cells.Item(136, 9).Value = 4786.643
$ws.Cells.Item(136, 10).Value = 9999
$ws.Cells.Item(136, 11).Value = 14359.929
$ws.Cells.Item(136, 12).Value = 29997
$ws.Cells.Item(136, 13).Value = -11809.929
$ws.Cells.Item(136, 14).Value = -35097

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 234.16667
$ws.Cells.Item(4, 10).Value = 351.25
$ws.Cells.Item(4, 12).Value = 351.25
$ws.Cells.Item(4, 14).Value = -581.25
$ws.Cells.Item(20, 8).Value = 10400
$ws.Cells.Item(20, 9).Value = 13390.1
$ws.Cells.Item(20, 10).Value = 2924.75
$ws.Cells.Item(20, 11).Value = 13390.1
$ws.Cells.Item(20, 12).Value = 2924.75
$ws.Cells.Item(20, 13).Value = -13143.1
$ws.Cells.Item(20, 14).Value = -3418.75
$ws.Cells.Item(22, 8).Value = 16257.143
$ws.Cells.Item(22, 9).Value = 25450
$ws.Cells.Item(22, 11).Value = 25450
$ws.Cells.Item(22, 13).Value = -25277
$ws.Cells.Item(39, 8).Value = 2800
$ws.Cells.Item(39, 10).Value = 2800
$ws.Cells.Item(39, 12).Value = 2800
$ws.Cells.Item(39, 14).Value = -3578

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1813.6
$ws.Cells.Item(16, 9).Value = 1034.5
$ws.Cells.Item(16, 10).Value = 2333
$ws.Cells.Item(16, 11).Value = 1034.5
$ws.Cells.Item(16, 12).Value = 2333
$ws.Cells.Item(16, 13).Value = -747.5
$ws.Cells.Item(16, 14).Value = -2907
$ws.Cells.Item(31, 8).Value = 4052.889
$ws.Cells.Item(31, 9).Value = 1772.5
$ws.Cells.Item(31, 10).Value = 8613.666999999999
$ws.Cells.Item(31, 11).Value = 1772.5
$ws.Cells.Item(31, 12).Value = 8613.666999999999
$ws.Cells.Item(31, 13).Value = -1477.5
$ws.Cells.Item(31, 14).Value = -9203.666999999999
$ws.Cells.Item(34, 8).Value = 4052.889
$ws.Cells.Item(34, 9).Value = 1772.5
$ws.Cells.Item(34, 10).Value = 8613.666999999999
$ws.Cells.Item(34, 11).Value = 1772.5
$ws.Cells.Item(34, 12).Value = 8613.666999999999
$ws.Cells.Item(34, 13).Value = -1570.5
$ws.Cells.Item(34, 14).Value = -9017.666999999999
$ws.Cells.Item(62, 8).Value = 4011
$ws.Cells.Item(62, 10).Value = 3908.5
$ws.Cells.Item(62, 12).Value = 3908.5
$ws.Cells.Item(62, 14).Value = -5156.5
$ws.Cells.Item(65, 8).Value = 4011
$ws.Cells.Item(65, 10).Value = 3908.5
$ws.Cells.Item(65, 12).Value = 19542.5
$ws.Cells.Item(65, 14).Value = -25782.5
$ws.Cells.Item(94, 8).Value = 1148.3846
$ws.Cells.Item(94, 10).Value = 1710.5714
$ws.Cells.Item(94, 12).Value = 1710.5714
$ws.Cells.Item(94, 14).Value = -2612.5714
$ws.Cells.Item(113, 8).Value = 1813.6
$ws.Cells.Item(113, 9).Value = 1034.5
$ws.Cells.Item(113, 10).Value = 2333
$ws.Cells.Item(113, 11).Value = 1034.5
$ws.Cells.Item(113, 12).Value = 2333
$ws.Cells.Item(113, 13).Value = 1135.5
$ws.Cells.Item(113, 14).Value = -6673

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 204.1
$ws.Cells.Item(2, 9).Value = 26.333334
$ws.Cells.Item(2, 10).Value = 280.2857
$ws.Cells.Item(2, 11).Value = 158.000004
$ws.Cells.Item(2, 12).Value = 1681.7142
$ws.Cells.Item(2, 13).Value = -45.00000399999999
$ws.Cells.Item(2, 14).Value = -1907.7142
$ws.Cells.Item(38, 8).Value = 20834000
$ws.Cells.Item(38, 9).Value = 21.166666
$ws.Cells.Item(38, 10).Value = 33334386
$ws.Cells.Item(38, 11).Value = 63.499998
$ws.Cells.Item(38, 12).Value = 100003158
$ws.Cells.Item(38, 13).Value = 283.500002
$ws.Cells.Item(38, 14).Value = -100003852
$ws.Cells.Item(109, 8).Value = 6263.5
$ws.Cells.Item(109, 9).Value = 916.2
$ws.Cells.Item(109, 11).Value = 2748.6
$ws.Cells.Item(109, 13).Value = -1708.6
$ws.Cells.Item(119, 8).Value = 964.5
$ws.Cells.Item(119, 9).Value = 964.5
$ws.Cells.Item(119, 11).Value = 2893.5
$ws.Cells.Item(119, 13).Value = 1944.5

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 8397.5
$ws.Cells.Item(2, 10).Value = 10.5
$ws.Cells.Item(2, 12).Value = 10.5
$ws.Cells.Item(2, 14).Value = -236.5
$ws.Cells.Item(70, 8).Value = 15891367
$ws.Cells.Item(70, 9).Value = 22238948
$ws.Cells.Item(70, 10).Value = 22414.666
$ws.Cells.Item(70, 11).Value = 22238948
$ws.Cells.Item(70, 12).Value = 22414.666
$ws.Cells.Item(70, 13).Value = -22238678
$ws.Cells.Item(70, 14).Value = -22954.666
$ws.Cells.Item(73, 8).Value = 15891367
$ws.Cells.Item(73, 9).Value = 22238948
$ws.Cells.Item(73, 10).Value = 22414.666
$ws.Cells.Item(73, 11).Value = 22238948
$ws.Cells.Item(73, 12).Value = 22414.666
$ws.Cells.Item(73, 13).Value = -22238012
$ws.Cells.Item(73, 14).Value = -24286.666
$ws.Cells.Item(126, 8).Value = 4099
$ws.Cells.Item(126, 9).Value = 3237.3333
$ws.Cells.Item(126, 10).Value = 4745.25
$ws.Cells.Item(126, 11).Value = 9711.999899999999
$ws.Cells.Item(126, 12).Value = 14235.75
$ws.Cells.Item(126, 13).Value = -7241.999899999999
$ws.Cells.Item(126, 14).Value = -19175.75
$ws.Cells.Item(132, 8).Value = 5999.3335
$ws.Cells.Item(132, 9).Value = 3199.4
$ws.Cells.Item(132, 10).Value = 19999
$ws.Cells.Item(132, 11).Value = 9598.200000000001
$ws.Cells.Item(132, 12).Value = 59997
$ws.Cells.Item(132, 14).Value = -65057
$ws.Cells.Item(132, 13).Value = -7068.200000000001

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3952.7144
$ws.Cells.Item(16, 9).Value = 665.85
$ws.Cells.Item(16, 10).Value = 69690
$ws.Cells.Item(16, 11).Value = 665.85
$ws.Cells.Item(16, 12).Value = 69690
$ws.Cells.Item(16, 13).Value = -495.85
$ws.Cells.Item(16, 14).Value = -70030
$ws.Cells.Item(61, 8).Value = 5213.7144
$ws.Cells.Item(61, 9).Value = 3124
$ws.Cells.Item(61, 11).Value = 3124
$ws.Cells.Item(61, 13).Value = -2922
$ws.Cells.Item(113, 8).Value = 5213.7144
$ws.Cells.Item(113, 9).Value = 3124
$ws.Cells.Item(113, 11).Value = 3124
$ws.Cells.Item(113, 13).Value = -954
$ws.Cells.Item(122, 8).Value = 3666.8928
$ws.Cells.Item(122, 9).Value = 3506.6538
$ws.Cells.Item(122, 10).Value = 5750
$ws.Cells.Item(122, 11).Value = 10519.9614
$ws.Cells.Item(122, 12).Value = 17250
$ws.Cells.Item(122, 13).Value = -8069.9614
$ws.Cells.Item(122, 14).Value = -22150

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(63, 8).Value = 19367.75
$ws.Cells.Item(63, 9).Value = 18226
$ws.Cells.Item(63, 11).Value = 18226
$ws.Cells.Item(63, 13).Value = -17602
$ws.Cells.Item(66, 8).Value = 19367.75
$ws.Cells.Item(66, 9).Value = 18226
$ws.Cells.Item(66, 11).Value = 54678
$ws.Cells.Item(66, 13).Value = -51558
$ws.Cells.Item(136, 8).Value = 4774.3657
$ws.Cells.Item(136, 9).Value = 3356.3447
$ws.Cells.Item(136, 10).Value = 8201.25
$ws.Cells.Item(136, 11).Value = 10069.0341
$ws.Cells.Item(136, 12).Value = 24603.75
$ws.Cells.Item(136, 13).Value = -7519.034100000001
$ws.Cells.Item(136, 14).Value = -29703.75
